# Update "想去人数" (F column) values on the sheets that hold the
# conference data: "展览" (sheet1) and "全部类型" (sheet4).
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1476
    "F3"  = 957
    "F5"  = 2216
    "F7"  = 1360
    "F9"  = 143
    "F11" = 324
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
